$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 13 data - reuse existing style (same as row 12) by copying format first
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A13").Value = 45747
$ws.Range("B13").Value = 0.375
$ws.Range("C13").Value = "Futconnect 3103 0900"
$ws.Range("D13").Value = "Arrumado estatisticas do jogo quando jogador cancela participação"
